# The diff touches the "Table 3" graphicFrame on slide 3: the word
# "CarsForSale" in the header cell (row 1, col 1) was re-typed starting
# at its second character ("C" | "arsForSale"), which is exactly what
# PowerPoint does when a user clicks into existing text, selects from
# the 2nd character to the end, and retypes it (the freshly-typed part
# becomes its own run while the untouched leading "C" keeps its
# original run/formatting).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shp = $s.Shapes.Item(3)          # "Table 3"
$tbl = $shp.Table

$cell = $tbl.Cell(1, 1)
$tr = $cell.Shape.TextFrame.TextRange

$fullText = $tr.Text
$len = $fullText.Length

# Leave the first character ("C") untouched and re-type everything
# from the 2nd character to the end ("arsForSale"), mirroring the
# select+retype edit recorded in the diff.
$retyped = $tr.Characters(2, $len - 1)
$retyped.Text = "arsForSale"

Write-Host "Table 3 / Cell(1,1):" $tr.Text
